$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '65.863.13'
$r.Style = "Normal"
$ws.Range("E2").Value = '  +1.54%  '

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '2.694.66'
$r.Style = "Normal"
$ws.Range("E3").Value = '  +2.19%  '

$ws.Range("E4").Value = '  +0.04%  '

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '608.10'
$r.Style = "Normal"
$ws.Range("E5").Value = '  +2.06%  '

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '157.93'
$r.Style = "Normal"
$ws.Range("E6").Value = '  +1.64%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  -0.61%  '

$ws.Range("E9").Value = '  +6.03%  '

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '6.01'
$r.Style = "Normal"
$ws.Range("E10").Value = '  +3.94%  '

$ws.Range("E11").Value = '  +1.14%  '

$ws.Range("E12").Value = '  +1.14%  '

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '30.35'
$r.Style = "Normal"
$ws.Range("E13").Value = '  +4.50%  '

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '0.0000200'
$r.Style = "Normal"
$ws.Range("E14").Value = '  +7.20%  '

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '3.180.89'
$r.Style = "Normal"
$ws.Range("E15").Value = '  +2.19%  '

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '65.719.31'
$r.Style = "Normal"
$ws.Range("E16").Value = '  +1.41%  '

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '2.700.22'
$r.Style = "Normal"
$ws.Range("E17").Value = '  +4.05%  '

$ws.Range("E19").Value = '  +2.14%  '

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '359.31'
$r.Style = "Normal"
$ws.Range("E20").Value = '  +2.25%  '

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '7.55'
$r.Style = "Normal"
$ws.Range("E21").Value = '  +3.72%  '

$ws.Range("E22").Value = '  -0.21%  '

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '70.68'
$r.Style = "Normal"
$ws.Range("E23").Value = '  +4.04%  '

$ws.Range("E24").Value = '  +3.69%  '

$ws.Range("E25").Value = '  +14.10%  '

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '1.67'
$r.Style = "Normal"
$ws.Range("E26").Value = '  -1.07%  '

$ws.Range("E27").Value = '  +3.08%  '

$ws.Range("E28").Value = '  +5.15%  '

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '8.41'
$r.Style = "Normal"
$ws.Range("E29").Value = '  +3.84%  '

$ws.Range("E30").Value = '  +5.05%  '

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '540.98'
$r.Style = "Normal"
$ws.Range("E31").Value = '  +6.02%  '

$ws.Range("E32").Value = '  +0.16%  '

$ws.Range("E33").Value = '  +2.09%  '

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '6.67'
$r.Style = "Normal"
$ws.Range("E34").Value = '  +6.01%  '

$ws.Range("E35").Value = '  -3.54%  '

$ws.Range("E36").Value = '  +1.89%  '

$ws.Range("E37").Value = '  +3.39%  '

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '163.37'
$r.Style = "Normal"
$ws.Range("E38").Value = '  -0.38%  '

$ws.Range("E39").Value = '  -0.02%  '

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E40").Value = '  +0.14%  '

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '171.67'
$r.Style = "Normal"
$ws.Range("E41").Value = '  +3.97%  '

$ws.Range("E42").Value = '  -0.01%  '

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '42.66'
$r.Style = "Normal"
$ws.Range("E43").Value = '  +1.15%  '

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '4.18'
$r.Style = "Normal"
$ws.Range("E44").Value = '  +2.74%  '

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.0615'
$r.Style = "Normal"
$ws.Range("E45").Value = '  +0.60%  '

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '23.54'
$r.Style = "Normal"
$ws.Range("E46").Value = '  +2.83%  '

$ws.Range("E47").Value = '  +4.37%  '

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '0.0266'
$r.Style = "Normal"
$ws.Range("E48").Value = '  +4.52%  '

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '0.657'
$r.Style = "Normal"
$ws.Range("E49").Value = '  +1.62%  '

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '21.01'
$r.Style = "Normal"
$ws.Range("E50").Value = '  +9.04%  '

$ws.Range("E51").Value = '  +1.31%  '
